$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 820259.9
$ws.Range("I15").Value = 820259.9
$ws.Range("K15").Value = 2460779.7
$ws.Range("M15").Value = -2460610.7
$ws.Range("H28").Value = 424.2963
$ws.Range("I28").Value = 297.8095
$ws.Range("J28").Value = 867
$ws.Range("K28").Value = 297.8095
$ws.Range("L28").Value = 867
$ws.Range("M28").Value = 187.1905
$ws.Range("N28").Value = -1837
$ws.Range("H116").Value = 2839.6667
$ws.Range("I116").Value = 2741.3333
$ws.Range("J116").Value = 3183.8333
$ws.Range("K116").Value = 2741.3333
$ws.Range("L116").Value = 3183.8333
$ws.Range("M116").Value = 700.6667000000002
$ws.Range("N116").Value = -10067.8333

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 25004134
$ws.Range("I32").Value = 30305322
$ws.Range("J32").Value = 12816.143
$ws.Range("K32").Value = 30305322
$ws.Range("L32").Value = 12816.143
$ws.Range("M32").Value = -30305035
$ws.Range("N32").Value = -13390.143
$ws.Range("H45").Value = 2873.9
$ws.Range("I45").Value = 2625
$ws.Range("J45").Value = 3039.8333
$ws.Range("K45").Value = 2625
$ws.Range("L45").Value = 3039.8333
$ws.Range("M45").Value = -2248
$ws.Range("N45").Value = -3793.8333
$ws.Range("H61").Value = 3772.1428
$ws.Range("I61").Value = 3681
$ws.Range("J61").Value = 4000
$ws.Range("K61").Value = 3681
$ws.Range("L61").Value = 4000
$ws.Range("M61").Value = -3469
$ws.Range("N61").Value = -4424
$ws.Range("H74").Value = 5620.4346
$ws.Range("I74").Value = 6774.278
$ws.Range("J74").Value = 1466.6
$ws.Range("K74").Value = 6774.278
$ws.Range("L74").Value = 1466.6
$ws.Range("M74").Value = -5900.278
$ws.Range("N74").Value = -3214.6
$ws.Range("H77").Value = 5620.4346
$ws.Range("I77").Value = 6774.278
$ws.Range("J77").Value = 1466.6
$ws.Range("K77").Value = 33871.39
$ws.Range("L77").Value = 7333
$ws.Range("M77").Value = -29503.39
$ws.Range("N77").Value = -16069
$ws.Range("H97").Value = 610
$ws.Range("I97").Value = 635.4545000000001
$ws.Range("J97").Value = 575
$ws.Range("K97").Value = 635.4545000000001
$ws.Range("L97").Value = 575
$ws.Range("M97").Value = -139.4545000000001
$ws.Range("N97").Value = -1567
$ws.Range("H132").Value = 2487.4285
$ws.Range("I132").Value = 2294.2307
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 6882.6921
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = -4352.6921
$ws.Range("N132").Value = -20057
$ws.Range("H136").Value = 3772.1428
$ws.Range("I136").Value = 3681
$ws.Range("J136").Value = 4000
$ws.Range("K136").Value = 11043
$ws.Range("L136").Value = 12000
$ws.Range("M136").Value = -8493
$ws.Range("N136").Value = -17100

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H92").Value = 46537.625
$ws.Range("J92").Value = 46537.625
$ws.Range("L92").Value = 46537.625
$ws.Range("N92").Value = -51529.625
$ws.Range("H134").Value = 1638.4857
$ws.Range("I134").Value = 1638.4857
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 4915.4571
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -2380.4571
$ws.Range("N134").ClearContents()

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1266.1082
$ws.Range("I58").Value = 552.619
$ws.Range("J58").Value = 2202.5625
$ws.Range("K58").Value = 552.619
$ws.Range("L58").Value = 2202.5625
$ws.Range("M58").Value = -349.619
$ws.Range("N58").Value = -2608.5625
$ws.Range("H117").Value = 47141.332
$ws.Range("J117").Value = 47141.332
$ws.Range("L117").Value = 47141.332
$ws.Range("N117").Value = -56319.332
$ws.Range("H132").Value = 1715.2727
$ws.Range("I132").Value = 1286.9
$ws.Range("J132").Value = 5999
$ws.Range("K132").Value = 3860.7
$ws.Range("L132").Value = 17997
$ws.Range("M132").Value = -1330.7
$ws.Range("N132").Value = -23057
$ws.Range("H134").Value = 3408.2917
$ws.Range("I134").Value = 1900.0454
$ws.Range("K134").Value = 5700.1362
$ws.Range("M134").Value = -3165.1362
$ws.Range("H136").Value = 1266.1082
$ws.Range("I136").Value = 552.619
$ws.Range("J136").Value = 2202.5625
$ws.Range("K136").Value = 1657.857
$ws.Range("L136").Value = 6607.6875
$ws.Range("M136").Value = 892.143
$ws.Range("N136").Value = -11707.6875

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 697.95
$ws.Range("I122").Value = 590
$ws.Range("J122").Value = 1021.8
$ws.Range("K122").Value = 5310
$ws.Range("L122").Value = 9196.199999999999
$ws.Range("M122").Value = -2860
$ws.Range("N122").Value = -14096.2
$ws.Range("H123").Value = 4728.5713
$ws.Range("I123").Value = 3600
$ws.Range("J123").Value = 4916.6665
$ws.Range("K123").Value = 10800
$ws.Range("L123").Value = 14749.9995
$ws.Range("M123").Value = -8350
$ws.Range("N123").Value = -19649.9995
$ws.Range("H137").Value = 5111.905
$ws.Range("I137").Value = 1668.75
$ws.Range("J137").Value = 7230.769
$ws.Range("K137").Value = 5006.25
$ws.Range("L137").Value = 21692.307
$ws.Range("M137").Value = 93.75
$ws.Range("N137").Value = -31892.307

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 9416.571
$ws.Range("I113").Value = 1986
$ws.Range("K113").Value = 1986
$ws.Range("M113").Value = 184
$ws.Range("H120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("N120").ClearContents()
$ws.Range("H132").Value = 3307.4285
$ws.Range("I132").Value = 2675.647
$ws.Range("J132").Value = 5992.5
$ws.Range("K132").Value = 8026.941
$ws.Range("L132").Value = 17977.5
$ws.Range("M132").Value = -5496.941
$ws.Range("N132").Value = -23037.5
$ws.Range("H141").Value = 69400
$ws.Range("J141").Value = 69400
$ws.Range("L141").Value = 69400
$ws.Range("N141").Value = -79760

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 3000
$ws.Range("I51").Value = 3000
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 3000
$ws.Range("L51").Value = 0
$ws.Range("M51").Value = -2490
$ws.Range("N51").ClearContents()
$ws.Range("H107").Value = 737.56525
$ws.Range("I107").Value = 657.4666999999999
$ws.Range("J107").Value = 887.75
$ws.Range("K107").Value = 1972.4001
$ws.Range("L107").Value = 2663.25
$ws.Range("M107").Value = -52.40009999999984
$ws.Range("N107").Value = -6503.25
$ws.Range("H136").Value = 1725.1471
$ws.Range("I136").Value = 1342.44
$ws.Range("J136").Value = 2788.2222
$ws.Range("K136").Value = 4027.32
$ws.Range("L136").Value = 8364.6666
$ws.Range("M136").Value = -1477.32
$ws.Range("N136").Value = -13464.6666
